$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 / J1 (style matches the rest of the header row: s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Cells.Item(1, 8).Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for I2:J86
$iValues = @(7,7,6,8,8,7,8,8,7,8,7,8,8,7,8,7,8,7,6,8,8,7,8,7,7,7,8,9,8,6,7,8,8,8,8,7,8,9,7,8,6,8,7,8,8,8,9,8,8,8,8,8,7,8,9,8,8,8,9,10,9,7,8,8,9,9,9,7,9,11,7,8,8,8,8,9,8,8,9,6,6,4,9,8,7)
$jValues = @(8,7,7,8,8,8,8,8,7,8,7,8,8,8,8,7,8,7,6,8,8,7,8,7,7,7,8,9,8,7,7,8,8,8,8,7,8,9,8,8,7,8,7,8,8,8,9,8,8,8,8,8,7,8,9,8,8,8,9,10,9,7,8,8,9,9,9,8,9,11,7,8,8,8,8,9,8,8,9,6,7,4,9,8,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
